$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 9-10 content down to rows 10-11, add two new rows (Introvid, Posts)
# at 12-13, and clear row 9 (keeping its original formatting).

# Preserve row10's content ("Statistic" task) into row11 before it gets overwritten.
$ws.Range("A10:E10").Copy($ws.Range("A11:E11"))

# Move row9's content ("Overview" task) into row10.
$ws.Range("A9:E9").Copy($ws.Range("A10:E10"))

# Row 9 becomes blank (formatting only, no values).
$ws.Range("A9:E9").ClearContents()

# Clone the formatting of row 11 into the two new rows 12 and 13.
$ws.Range("A11:E11").Copy($ws.Range("A12:E12"))
$ws.Range("A11:E11").Copy($ws.Range("A13:E13"))

# Fill in the new "Introvid" task on row 12.
$ws.Range("A12").Value = "Introvid"
$ws.Range("B12").Value = 45068
$ws.Range("C12").Value = 45068
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = "Nhật Kha"

# Fill in the new "Posts" task on row 13.
$ws.Range("A13").Value = "Posts"
$ws.Range("B13").Value = 45069
$ws.Range("C13").Value = 45069
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = "Nhật Kha"

# Widen column B to fit the date values (best-fit-like width).
$ws.Columns.Item(2).ColumnWidth = 9.8

# Update the view: scroll down a bit and select A16, like in the saved file.
$ws.Range("A16").Select()
